$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "62.649.69"
    "E2"  = "  +1.92%  "
    "D3"  = "2.933.72"
    "E3"  = "  +0.00%  "
    "E4"  = "  +0.16%  "
    "D5"  = "592.66"
    "E5"  = "  -0.73%  "
    "D6"  = "146.16"
    "E6"  = "  +0.57%  "
    "E7"  = "  +0.15%  "
    "D8"  = "2.931.66"
    "E8"  = "  -0.05%  "
    "E9"  = "  +0.53%  "
    "D10" = "7.27"
    "E10" = "  +3.69%  "
    "D11" = "0.149"
    "E11" = "  +4.86%  "
    "D12" = "0.439"
    "D13" = "0.0000235"
    "E13" = "  +4.20%  "
    "D14" = "32.51"
    "E14" = "  -3.22%  "
    "E15" = "  -0.86%  "
    "D16" = "3.419.84"
    "E16" = "  -0.03%  "
    "D17" = "62.670.85"
    "E17" = "  +2.00%  "
    "D18" = "6.65"
    "D19" = "2.932.41"
    "E19" = "  -0.08%  "
    "D20" = "437.53"
    "E20" = "  +1.28%  "
    "D21" = "13.30"
    "E21" = "  -1.30%  "
    "D22" = "0.662"
    "E22" = "  -2.18%  "
    "D23" = "7.00"
    "E23" = "  -1.41%  "
    "D24" = "80.76"
    "E24" = "  -1.52%  "
    "D25" = "11.02"
    "E25" = "  +0.80%  "
    "D26" = "2.11"
    "E26" = "  -3.73%  "
    "D27" = "11.66"
    "E27" = "  -1.05%  "
    "D28" = "0.999"
    "E28" = "  -0.06%  "
    "E29" = "  -0.07%  "
    "D30" = "7.15"
    "E30" = "  +3.07%  "
    "D31" = "2.60"
    "E31" = "  -0.48%  "
    "D32" = "0.0000100"
    "E32" = "  +13.05%  "
    "E33" = "  -1.26%  "
    "D34" = "26.21"
    "E34" = "  -1.68%  "
    "E35" = "  +0.06%  "
    "D36" = "0.986"
    "E36" = "  -2.56%  "
    "D37" = "3.07"
    "E37" = "  +2.51%  "
    "E38" = "  -1.71%  "
    "D39" = "49.63"
    "E39" = "  -0.79%  "
    "E40" = "  +0.26%  "
    "D41" = "8.43"
    "E41" = "  -1.78%  "
    "D42" = "0.116"
    "E42" = "  -5.40%  "
    "D43" = "0.278"
    "E43" = "  -1.52%  "
    "D44" = "38.79"
    "E44" = "  -8.35%  "
    "D45" = "2.694.16"
    "E45" = "  -0.11%  "
    "D46" = "134.20"
    "E46" = "  +0.01%  "
    "D47" = "360.01"
    "E47" = "  -0.61%  "
    "E48" = "  -3.68%  "
    "E49" = "  -0.01%  "
    "D51" = "22.53"
    "E51" = "  -4.94%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
